$d = $word.ActiveDocument

# Locate the target paragraphs by their distinctive content rather than by
# hard-coded indices, so the script is robust to how the paragraph
# collection is walked.
$emptyIdx = New-Object System.Collections.ArrayList
$longRelIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.Length -eq 1) {
        [void]$emptyIdx.Add($i)
    }
    if ($t.Length -gt 10 -and $t.Substring(0,1) -eq [char]0x5728) {
        # starts with "在" -> "在软件工程领域，可以盲目的认为关系按照对应数量，分为以下3中关系。..."
        $longRelIdx = $i
    }
}

# --- Hunk 1a: the first empty paragraph (right after "...或“对象”。") loses
# its stray <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>,
# becoming a completely bare paragraph (<w:p/>).
$firstEmpty = $emptyIdx[0]
$p1 = $d.Paragraphs($firstEmpty)
$p1.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")

# --- Hunk 1b: the following centered paragraph ("工程师是程序世界的神，做神很难")
# keeps its centering but loses the stray rFonts hint on the paragraph mark.
$p2 = $d.Paragraphs($firstEmpty + 1)
$p2.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='1EA0E6AE' w14:textId='6B815DBA' w:rsidR='00BA0593' w:rsidRPr='00631EA0' w:rsidRDefault='00381201' w:rsidP='00BA0593'><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>工程师</w:t></w:r><w:r w:rsidR='00BA0593'><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>是程序世界的神，做神很难</w:t></w:r></w:p>")

# --- Hunk 2: split the run containing "...分为以下3中关系。分为..." into three
# runs, carving out the lone "种" character; same run formatting throughout.
$p3 = $d.Paragraphs($longRelIdx)
$p3.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='0195D100' w14:textId='29B0A21E' w:rsidR='0073195D' w:rsidRPr='00631EA0' w:rsidRDefault='0073195D' w:rsidP='00631EA0'><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>在软件工程领域，可以盲目的认为关系按照对应数量，分为以下3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>种</w:t></w:r><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>关系。分为“对象”和“类”两个级别描述。</w:t></w:r></w:p>")

# --- Hunk 3: the second empty paragraph (right after "轮询") loses its
# stray <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>,
# becoming a completely bare paragraph (<w:p/>), same pattern as hunk 1a.
$secondEmpty = $emptyIdx[1]
$p4 = $d.Paragraphs($secondEmpty)
$p4.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")
